$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("A3").Value = "Third line"

$ws.Range("A3").Select()
